# Update the VytrackUsers sheet so every generated test-user row shares the
# same placeholder first/last name ("John"/"Doe") instead of the random
# fake names that used to live in columns C and D (rows 3-7).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3").Value = "John"
$ws.Range("D3").Value = "Doe"

$ws.Range("C4").Value = "John"
$ws.Range("D4").Value = "Doe"

$ws.Range("C5").Value = "John"
$ws.Range("D5").Value = "Doe"

$ws.Range("C6").Value = "John"
$ws.Range("D6").Value = "Doe"

$ws.Range("C7").Value = "John"
$ws.Range("D7").Value = "Doe"

# Match the saved selection / zoom state recorded in the updated workbook.
$ws.Range("L14").Select() | Out-Null
$excel.ActiveWindow.Zoom = 100
